# "Added transformation and outliers"
#
# The Jan 10 ("M") class session is split into a Fri/Mon combo session
# (row 6, now "F/M") and its own transformations session (row 7), and the
# class that used to be taught on row 8 moves from Monday to Tuesday
# (C8 "M" -> "T", D8 date +1 day) and is retitled from "Multiple Linear
# Regression" to "Intro to Multiple Linear Regression". The new
# transformations session (row 7) gets slides/ae/hw links, and the MLR
# session (row 8) gets a project-instructions link.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: day-of-week label becomes a combined Friday/Monday slot.
$ws.Range("C6").Value = "F/M"

# Row 7: fill in the new Transformations session's resource links.
# Clone the formatting already used by a populated cell in this row so the
# new cells pick up the same style as their siblings (I7 should look like
# the other "data link" cells, e.g. J7).
$ws.Range("J7").Copy()
$ws.Range("I7").PasteSpecial(-4122)
$ws.Range("G7").Value = "/slides/06-slr-transformations.qmd"
$ws.Range("H7").Value = "/ae/ae-06-transformations.qmd"
$ws.Range("I7").Value = "/hw/hw-04.qmd"

# Row 8: Multiple Linear Regression intro moves from Monday to Tuesday.
$ws.Range("F7").Copy()
$ws.Range("J8").PasteSpecial(-4122)
$ws.Range("C8").Value = "T"
$ws.Range("D8").Value = 45668
$ws.Range("J8").Value = "/project/project-instructions.qmd#plots-tables"
$ws.Range("E8").Value = "Intro to Multiple Linear Regression"

# Restore the active cell/selection to E9 (matches the saved view state).
$ws.Range("E9").Select()
